# Updates recomputed profit figures across the Kraken_Profits sheets (ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR). Values in columns H-N of specific rows are
# refreshed to reflect newly recalculated cost/profit numbers; a few rows gain
# new M/N totals while a few others have their (now-empty) M/N cells cleared.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1690.6
$ws.Range("I15").Value = 1690.6
$ws.Range("K15").Value = 5071.799999999999
$ws.Range("M15").Value = -4902.799999999999
$ws.Range("H19").Value = 467.72726
$ws.Range("I19").Value = 585.8333
$ws.Range("J19").Value = 326
$ws.Range("K19").Value = 585.8333
$ws.Range("L19").Value = 326
$ws.Range("M19").Value = -410.8333
$ws.Range("N19").Value = -676
$ws.Range("H28").Value = 3420.8333
$ws.Range("I28").Value = 3368.182
$ws.Range("K28").Value = 3368.182
$ws.Range("M28").Value = -2883.182
$ws.Range("H40").Value = 6133.75
$ws.Range("I40").Value = 2267.5
$ws.Range("K40").Value = 2267.5
$ws.Range("M40").Value = -2092.5
$ws.Range("H64").Value = 3999
$ws.Range("I64").Value = 3999
$ws.Range("K64").Value = 3999
$ws.Range("M64").Value = -3751
$ws.Range("H67").Value = 3999
$ws.Range("I67").Value = 3999
$ws.Range("K67").Value = 3999
$ws.Range("M67").Value = -3141
$ws.Range("H107").Value = 1576.4
$ws.Range("I107").Value = 1024.1428
$ws.Range("K107").Value = 1024.1428
$ws.Range("M107").Value = 895.8571999999999
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 1754
$ws.Range("H116").Value = 4987
$ws.Range("I116").Value = 4987
$ws.Range("J116").Value = 4987
$ws.Range("K116").Value = 4987
$ws.Range("L116").Value = 4987
$ws.Range("M116").Value = -1545
$ws.Range("N116").Value = -11871

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 504.4
$ws.Range("I2").Value = 504.4
$ws.Range("K2").Value = 504.4
$ws.Range("M2").Value = -391.4
$ws.Range("H45").Value = 2371.4243
$ws.Range("I45").Value = 2171.6538
$ws.Range("J45").Value = 3113.4285
$ws.Range("K45").Value = 2171.6538
$ws.Range("L45").Value = 3113.4285
$ws.Range("M45").Value = -1794.6538
$ws.Range("N45").Value = -3867.4285
$ws.Range("J110").Value = 7000
$ws.Range("L110").Value = 7000
$ws.Range("N110").Value = -11090
$ws.Range("H116").Value = 504.4
$ws.Range("I116").Value = 504.4
$ws.Range("K116").Value = 504.4
$ws.Range("M116").Value = 1789.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 504.4
$ws.Range("I3").Value = 504.4
$ws.Range("K3").Value = 504.4
$ws.Range("M3").Value = -390.4
$ws.Range("H22").Value = 2899.5
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5346
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H80").Value = 830.2
$ws.Range("I80").Value = 748.6667
$ws.Range("J80").Value = 952.5
$ws.Range("K80").Value = 748.6667
$ws.Range("L80").Value = 952.5
$ws.Range("M80").Value = 249.3333
$ws.Range("N80").Value = -2948.5
$ws.Range("H83").Value = 830.2
$ws.Range("I83").Value = 748.6667
$ws.Range("J83").Value = 952.5
$ws.Range("K83").Value = 3743.3335
$ws.Range("L83").Value = 4762.5
$ws.Range("M83").Value = 1248.6665
$ws.Range("N83").Value = -14746.5
$ws.Range("H86").Value = 4530.769
$ws.Range("I86").Value = 3100
$ws.Range("J86").Value = 5425
$ws.Range("K86").Value = 3100
$ws.Range("L86").Value = 5425
$ws.Range("M86").Value = -1977
$ws.Range("N86").Value = -7671
$ws.Range("H89").Value = 4530.769
$ws.Range("I89").Value = 3100
$ws.Range("J89").Value = 5425
$ws.Range("K89").Value = 15500
$ws.Range("L89").Value = 27125
$ws.Range("M89").Value = -9884
$ws.Range("N89").Value = -38357
$ws.Range("H94").Value = 1425.8889
$ws.Range("I94").Value = 1404.7142
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 1404.7142
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -953.7141999999999
$ws.Range("N94").Value = -2402
$ws.Range("H105").Value = 2200
$ws.Range("I105").Value = 1400
$ws.Range("K105").Value = 1400
$ws.Range("M105").Value = 347

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H58").Value = 1838.7693
$ws.Range("I58").Value = 1658.6666
$ws.Range("K58").Value = 1658.6666
$ws.Range("M58").Value = -1455.6666
$ws.Range("H69").Value = 7000
$ws.Range("I69").Value = 7000
$ws.Range("K69").Value = 7000
$ws.Range("M69").Value = -6251
$ws.Range("H72").Value = 7000
$ws.Range("I72").Value = 7000
$ws.Range("K72").Value = 21000
$ws.Range("M72").Value = -17256
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H136").Value = 1838.7693
$ws.Range("I136").Value = 1658.6666
$ws.Range("K136").Value = 4975.9998
$ws.Range("M136").Value = -2425.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 25.866667
$ws.Range("J12").Value = 29.285715
$ws.Range("L12").Value = 87.857145
$ws.Range("N12").Value = -433.857145
$ws.Range("H23").Value = 2039.1
$ws.Range("I23").Value = 1899.3334
$ws.Range("J23").Value = 2248.75
$ws.Range("K23").Value = 5698.0002
$ws.Range("L23").Value = 6746.25
$ws.Range("M23").Value = -5463.0002
$ws.Range("N23").Value = -7216.25
$ws.Range("H80").Value = 4500.636
$ws.Range("I80").Value = 1626.75
$ws.Range("K80").Value = 4880.25
$ws.Range("M80").Value = -3944.25
$ws.Range("H83").Value = 4500.636
$ws.Range("I83").Value = 1626.75
$ws.Range("K83").Value = 14640.75
$ws.Range("M83").Value = -9960.75
$ws.Range("H113").Value = 578
$ws.Range("I113").Value = 578
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1734
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 436
$ws.Range("N113").ClearContents()
$ws.Range("H120").Value = 9000
$ws.Range("I120").Value = 9000
$ws.Range("K120").Value = 27000
$ws.Range("M120").Value = -22162
$ws.Range("H129").Value = 2347.5454
$ws.Range("I129").Value = 1403.5714
$ws.Range("J129").Value = 3999.5
$ws.Range("K129").Value = 4210.7142
$ws.Range("L129").Value = 11998.5
$ws.Range("M129").Value = 789.2857999999997
$ws.Range("N129").Value = -21998.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3999.5
$ws.Range("I70").Value = 3999.5
$ws.Range("K70").Value = 3999.5
$ws.Range("M70").Value = -3729.5
$ws.Range("H73").Value = 3999.5
$ws.Range("I73").Value = 3999.5
$ws.Range("K73").Value = 3999.5
$ws.Range("M73").Value = -3063.5
$ws.Range("H80").Value = 2794
$ws.Range("I80").Value = 2688
$ws.Range("K80").Value = 2688
$ws.Range("M80").Value = -1690
$ws.Range("H83").Value = 2794
$ws.Range("I83").Value = 2688
$ws.Range("K83").Value = 13440
$ws.Range("M83").Value = -8448
$ws.Range("H97").Value = 1239.7273
$ws.Range("I97").Value = 867.3333
$ws.Range("J97").Value = 1686.6
$ws.Range("K97").Value = 867.3333
$ws.Range("L97").Value = 1686.6
$ws.Range("M97").Value = -371.3333
$ws.Range("N97").Value = -2678.6
$ws.Range("H107").Value = 780.4
$ws.Range("I107").Value = 650.5
$ws.Range("K107").Value = 650.5
$ws.Range("M107").Value = 1269.5
$ws.Range("H136").Value = 80000
$ws.Range("J136").Value = 80000
$ws.Range("L136").Value = 240000
$ws.Range("N136").Value = -245100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 202.5
$ws.Range("I16").Value = 202.5
$ws.Range("K16").Value = 202.5
$ws.Range("M16").Value = -32.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H130").Value = 12485.6
$ws.Range("J130").Value = 12485.6
$ws.Range("L130").Value = 12485.6
$ws.Range("N130").Value = -22525.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 958.6875
$ws.Range("I107").Value = 568.25
$ws.Range("J107").Value = 1349.125
$ws.Range("K107").Value = 1704.75
$ws.Range("L107").Value = 4047.375
$ws.Range("M107").Value = 215.25
$ws.Range("N107").Value = -7887.375
$ws.Range("H126").Value = 681
$ws.Range("I126").Value = 714
$ws.Range("K126").Value = 2142
$ws.Range("M126").Value = 328
